# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume/coin-name updates described by the diff
# (commit: "Updated cryptos list on Thu Nov 16 16:23:13 UTC 2023 with GitHub Actions").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (safe to assign directly; Excel will not reinterpret them as numbers)
$textUpdates = @(
    @{ Cell = 'D2'; Value = '36.637.51' }
    @{ Cell = 'E2'; Value = '  +0.71%  ' }
    @{ Cell = 'D3'; Value = '2.007.29' }
    @{ Cell = 'E3'; Value = '  -0.41%  ' }
    @{ Cell = 'E4'; Value = '  +0.04%  ' }
    @{ Cell = 'E5'; Value = '  -1.65%  ' }
    @{ Cell = 'E6'; Value = '  -1.55%  ' }
    @{ Cell = 'E7'; Value = '  -0.15%  ' }
    @{ Cell = 'E8'; Value = '  -0.03%  ' }
    @{ Cell = 'E9'; Value = '  +3.33%  ' }
    @{ Cell = 'E10'; Value = '  -2.21%  ' }
    @{ Cell = 'E11'; Value = '  +4.61%  ' }
    @{ Cell = 'E12'; Value = '  -0.36%  ' }
    @{ Cell = 'B13'; Value = 'Polygon' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = 'E13'; Value = '  -1.75%  ' }
    @{ Cell = 'B14'; Value = 'Avalanche' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Cell = 'E14'; Value = '  +13.51%  ' }
    @{ Cell = 'E15'; Value = '  -3.52%  ' }
    @{ Cell = 'D16'; Value = '2.300.22' }
    @{ Cell = 'E16'; Value = '  -0.46%  ' }
    @{ Cell = 'E17'; Value = '  +1.83%  ' }
    @{ Cell = 'D18'; Value = '2.006.58' }
    @{ Cell = 'E18'; Value = '  -0.33%  ' }
    @{ Cell = 'D19'; Value = '36.529.30' }
    @{ Cell = 'E19'; Value = '  +0.54%  ' }
    @{ Cell = 'E20'; Value = '  -0.15%  ' }
    @{ Cell = 'D21'; Value = '0.0₃0875' }
    @{ Cell = 'E21'; Value = '  +1.55%  ' }
    @{ Cell = 'E22'; Value = '  +0.34%  ' }
    @{ Cell = 'E23'; Value = '  +0.74%  ' }
    @{ Cell = 'E24'; Value = '  +0.10%  ' }
    @{ Cell = 'E25'; Value = '  -6.07%  ' }
    @{ Cell = 'E26'; Value = '  -0.25%  ' }
    @{ Cell = 'E27'; Value = '  +3.10%  ' }
    @{ Cell = 'E28'; Value = '  +22.81%  ' }
    @{ Cell = 'E29'; Value = '  -2.06%  ' }
    @{ Cell = 'E30'; Value = '  +3.08%  ' }
    @{ Cell = 'E31'; Value = '  +0.13%  ' }
    @{ Cell = 'E32'; Value = '  -0.17%  ' }
    @{ Cell = 'E33'; Value = '  -2.18%  ' }
    @{ Cell = 'E34'; Value = '  +2.25%  ' }
    @{ Cell = 'E35'; Value = '  -2.40%  ' }
    @{ Cell = 'E36'; Value = '  +11.23%  ' }
    @{ Cell = 'E37'; Value = '  -4.55%  ' }
    @{ Cell = 'E38'; Value = '  +0.12%  ' }
    @{ Cell = 'E39'; Value = '  +0.80%  ' }
    @{ Cell = 'E40'; Value = '  +23.36%  ' }
    @{ Cell = 'B41'; Value = 'TrustWalletToken' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'E41'; Value = '  +3.41%  ' }
    @{ Cell = 'B42'; Value = 'Cronos' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'E42'; Value = '  -4.04%  ' }
    @{ Cell = 'E43'; Value = '  -0.11%  ' }
    @{ Cell = 'E44'; Value = '  -0.20%  ' }
    @{ Cell = 'E45'; Value = '  -0.66%  ' }
    @{ Cell = 'E46'; Value = '  -0.26%  ' }
    @{ Cell = 'E47'; Value = '  -1.41%  ' }
    @{ Cell = 'E48'; Value = '  -3.10%  ' }
    @{ Cell = 'D49'; Value = '1.359.72' }
    @{ Cell = 'E49'; Value = '  -5.94%  ' }
    @{ Cell = 'E50'; Value = '  -2.31%  ' }
    @{ Cell = 'D51'; Value = '2.192.05' }
    @{ Cell = 'E51'; Value = '  -0.39%  ' }
)

# Updates whose new value looks like a plain number (e.g. "247.26").
# The source cells are text cells (inline strings), so we force the destination
# cell to remain text too -- otherwise Excel's Value setter auto-coerces the
# string into a floating point number (losing formatting / introducing FP error).
$forceTextUpdates = @(
    @{ Cell = 'D5'; Value = '247.26' }
    @{ Cell = 'D6'; Value = '0.632' }
    @{ Cell = 'D7'; Value = '61.93' }
    @{ Cell = 'D9'; Value = '0.383' }
    @{ Cell = 'D10'; Value = '57.75' }
    @{ Cell = 'D11'; Value = '0.0778' }
    @{ Cell = 'D12'; Value = '0.104' }
    @{ Cell = 'D13'; Value = '0.888' }
    @{ Cell = 'D14'; Value = '22.97' }
    @{ Cell = 'D15'; Value = '14.28' }
    @{ Cell = 'D17'; Value = '5.53' }
    @{ Cell = 'D20'; Value = '71.93' }
    @{ Cell = 'D22'; Value = '5.32' }
    @{ Cell = 'D23'; Value = '236.22' }
    @{ Cell = 'D25'; Value = '2.51' }
    @{ Cell = 'D27'; Value = '9.89' }
    @{ Cell = 'D28'; Value = '0.137' }
    @{ Cell = 'D29'; Value = '159.82' }
    @{ Cell = 'D30'; Value = '20.21' }
    @{ Cell = 'D31'; Value = '0.121' }
    @{ Cell = 'D32'; Value = '1.18' }
    @{ Cell = 'D33'; Value = '4.99' }
    @{ Cell = 'D34'; Value = '0.0622' }
    @{ Cell = 'D35'; Value = '4.46' }
    @{ Cell = 'D36'; Value = '6.53' }
    @{ Cell = 'D37'; Value = '2.35' }
    @{ Cell = 'D40'; Value = '3.22' }
    @{ Cell = 'D41'; Value = '1.27' }
    @{ Cell = 'D42'; Value = '0.101' }
    @{ Cell = 'D44'; Value = '1.13' }
    @{ Cell = 'D45'; Value = '0.0215' }
    @{ Cell = 'D47'; Value = '93.36' }
    @{ Cell = 'D48'; Value = '7.61' }
    @{ Cell = 'D50'; Value = '2.88' }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

foreach ($u in $forceTextUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = '@'
    $cell.Value = $u.Value
    $cell.Style = 'Normal'
}

Write-Output ("Applied " + ($textUpdates.Count + $forceTextUpdates.Count) + " cell updates")
